$d = $word.ActiveDocument

# --- Change 1: remove the _GoBack bookmark currently sitting after
# "Public username" (it will be re-created at the end of the doc below,
# next to the last edit). ---
[void]$d.Bookmarks.Item("_GoBack").Delete()

# --- Change 2: add a new list item "Password unchanged for too long"
# right after the "Password is weak" bullet, in the same list
# (numId=1) and using the same paragraph/run formatting. ---
$weakIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Password is weak`r") {
        $weakIndex = $i
    }
}

$weakPara = $d.Paragraphs.Item($weakIndex)
$afterWeak = $weakPara.Range
$afterWeak.Collapse(0)
[void]$afterWeak.InsertParagraphAfter()

# Fill the freshly-created (empty) paragraph that now follows "Password is weak"
$newPara = $d.Paragraphs.Item($weakIndex + 1)
$fillRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' + `
  '<w:t xml:space="preserve">Password </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' + `
  '<w:t>unchanged</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' + `
  '<w:t xml:space="preserve"> for too long</w:t></w:r>' + `
  '</w:p>'
[void]$fillRange.InsertXML($newFrag)

# --- Change 3: split the "Information required for password reset
# available  publicly" sentence into separate runs, append
# " on another account", and move the _GoBack bookmark to the end of
# this (now last) paragraph. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$replaceRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$infoFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' + `
  '<w:t xml:space="preserve">Information required for password reset </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' + `
  '<w:t>available</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' + `
  '<w:t xml:space="preserve"> publicly</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' + `
  '<w:t xml:space="preserve"> on another account</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
[void]$replaceRange.InsertXML($infoFrag)
